$d = $word.ActiveDocument

function Set-HeadingNumber($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false,
                             $true, 1, $false, $newText, 2)
}

Set-HeadingNumber "Introduction" "1 Introduction"
Set-HeadingNumber "Why Use Quarto?" "2 Why Use Quarto?"
Set-HeadingNumber "Basic Example" "3 Basic Example"
Set-HeadingNumber "Quarto Features" "4 Quarto Features"
Set-HeadingNumber "Callouts" "4.1 Callouts"
Set-HeadingNumber "Code Chunks with Modern Syntax" "4.2 Code Chunks with Modern Syntax"
Set-HeadingNumber "Equations" "5 Equations"
Set-HeadingNumber "Tables" "6 Tables"
Set-HeadingNumber "Summary" "7 Summary"
Set-HeadingNumber "References" "8 References"
